# Small positional bug fixes for the pie-chart label textboxes (tx9..tx18)
# inside the group shape on slide 1. Only the positions (a:off x/y) move;
# sizes (a:ext) are unchanged.
#
# NOTE: PowerPoint COM reports/accepts shape Left/Top in points, while the
# underlying OOXML stores EMUs (1 pt = 12700 EMU). The literal point values
# below were chosen so that, after the host's internal float handling,
# they round-trip to the exact target EMU coordinates.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The 10 label textboxes all live inside the single group shape on the slide.
$grp = $s.Shapes.Item(2)
$items = $grp.GroupItems

$items.Item("tx9").Left  = 355.4623622047244
$items.Item("tx9").Top   = 204.40984351968504

$items.Item("tx10").Left = 383.4455118110236
$items.Item("tx10").Top  = 228.53457692913386

$items.Item("tx11").Left = 418.18818897637794
$items.Item("tx11").Top  = 260.32055118110236

$items.Item("tx12").Left = 464.6726871653543
$items.Item("tx12").Top  = 284.4452755905512

$items.Item("tx13").Left = 503.9407874015748
$items.Item("tx13").Top  = 309.15874015748034

$items.Item("tx14").Left = 505.3495375590551
$items.Item("tx14").Top  = 336.476062992126

$items.Item("tx15").Left = 326.6102462204724
$items.Item("tx15").Top  = 429.1151281102362

$items.Item("tx16").Left = 335.6170866141732
$items.Item("tx16").Top  = 456.43245094488185

$items.Item("tx17").Left = 319.9167029133858
$items.Item("tx17").Top  = 257.1011811023622

$items.Item("tx18").Left = 327.02708661417324
$items.Item("tx18").Top  = 284.41850393700787
